# "Actualizar" automation refresh: append a fresh 14-row disponibilidad
# block (rows 702-715) to Sheet1, mirroring the cyclical 14-row pattern
# already present in the sheet (Nombre / URL-hyperlink / Disponibilidad /
# Fecha). Also re-stamps the D688:D701 timestamps that the prior run wrote
# with a slightly different float literal for the same instant.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-stamp the previous batch's Fecha column (D688:D701) -----------
# Same instant, just re-serialized (the automation recomputed/rewrote the
# value when it appended the new batch below).
for ($r = 688; $r -le 701; $r++) {
    $ws.Range("D$r").Value = 44232.70380440972
}

# --- 2) Append the new batch: rows 702-715 --------------------------------
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$displayUrls = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
# The hyperlink's underlying Address (relationship Target); for MapStore
# the trailing "#/" lives in SubAddress instead, matching the workbook's
# existing rows (e.g. B10, B24, ...).
$linkAddresses = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$linkSubAddresses = @("","","","","","","","","/","","","","","")

$availability = "Disponible"
$fecha = 44232.72488631029

$startRow = 702
for ($i = 0; $i -lt 14; $i++) {
    $r = $startRow + $i

    $ws.Range("A$r").Value = $names[$i]
    $ws.Range("C$r").Value = $availability

    $ws.Range("D$r").Value = $fecha
    $ws.Range("D$r").NumberFormat = $ws.Range("D2").NumberFormat

    $ws.Range("B$r").Value = $displayUrls[$i]
    if ($linkSubAddresses[$i] -ne "") {
        $h = $ws.Hyperlinks.Add($ws.Range("B$r"), $linkAddresses[$i], $linkSubAddresses[$i])
    } else {
        $h = $ws.Hyperlinks.Add($ws.Range("B$r"), $linkAddresses[$i])
    }
    $ws.Range("B$r").Style = $ws.Range("B2").Style
}
